{"js": "// Load all paragraphs in the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// --- 1. Insert a new \"Meta description\" paragraph right after the\n//        first paragraph (the Heading1 title). ---\nconst titlePara = paragraphs.items[0];\nconst metaText =\n  \"Meta description: Experience the magic of the Cinderella fairy tale in this well-crafted slot game. Play Cinderella's Ball free for a chance to win big with interesting bonus rounds.\";\nconst metaPara = titlePara.insertParagraph(metaText, \"After\");\nmetaPara.style = \"Normal\";\nawait context.sync();\n\n// Make just the \"Meta description\" portion bold (separate run), leaving\n// the rest of the sentence in normal formatting.\nconst boldResults = metaPara.search(\"Meta description\", { matchCase: true });\nboldResults.load(\"items\");\nawait context.sync();\nboldResults.items[0].font.bold = true;\nawait context.sync();\n\n// This document's body paragraphs consistently start with a leading\n// empty run (matching the surrounding paragraphs' structure); reproduce\n// that here too.\nmetaPara.insertText(\"\", \"Start\");\nawait context.sync();\n\n// --- 2. Remove the duplicate bold title paragraph near the end, and\n//        3. replace the italic paragraph's text with the new prompt. ---\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs2.items;\nconst oldTitle = \"Play Cinderella's Ball Free | Exciting Fairy-tale Themed Slot Game\";\nconst oldDescription =\n  \"Experience the magic of the Cinderella fairy tale in this well-crafted slot game. Play Cinderella's Ball free for a chance to win big with interesting bonus rounds.\";\n\nlet dupTitlePara = null;\nlet descriptionPara = null;\n\nfor (let i = 1; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t === oldTitle) {\n    dupTitlePara = items[i];\n  } else if (t === oldDescription) {\n    descriptionPara = items[i];\n  }\n}\n\nif (dupTitlePara) {\n  dupTitlePara.delete();\n}\n\nif (descriptionPara) {\n  descriptionPara.insertText(\n    'Prompt: Create a cartoon-style feature image that features a happy Maya warrior wearing glasses. The image should be eye-catching and appealing, with bright colors and playful elements. The warrior should have a big smile on their face and be surrounded by symbols and elements from the Cinderella\\'s Ball game, such as the pumpkin carriage, the crystal slipper, and the magic wand. The text \"Cinderella\\'s Ball\" should be featured prominently in the image, in fun and playful font.',\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Insert a new \"Meta description\" paragraph right after the\n#        first paragraph (the Heading1 title). ---\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Style = \"Normal\"\n\n$metaRange = $metaPara.Range\n$metaRange.InsertAfter(\"Meta description: Experience the magic of the Cinderella fairy tale in this well-crafted slot game. Play Cinderella's Ball free for a chance to win big with interesting bonus rounds.\")\n\n# Bold just the \"Meta description\" portion (leaves the rest normal).\n$boldRange = $metaPara.Range\n$boldRange.Find.ClearFormatting()\n[void]$boldRange.Find.Execute(\"Meta description\")\n$boldRange.Font.Bold = 1\n\n# This document's body paragraphs consistently start with a leading\n# empty run (matching the surrounding paragraphs' structure); reproduce\n# that here too.\n$metaPara.Range.InsertBefore(\"\")\n\n# --- 2. Remove the duplicate bold title paragraph near the end, and\n#        3. replace the italic paragraph's text with the new prompt. ---\n$oldTitle = \"Play Cinderella's Ball Free | Exciting Fairy-tale Themed Slot Game\"\n$oldDescription = \"Experience the magic of the Cinderella fairy tale in this well-crafted slot game. Play Cinderella's Ball free for a chance to win big with interesting bonus rounds.\"\n$newPrompt = 'Prompt: Create a cartoon-style feature image that features a happy Maya warrior wearing glasses. The image should be eye-catching and appealing, with bright colors and playful elements. The warrior should have a big smile on their face and be surrounded by symbols and elements from the Cinderella''s Ball game, such as the pumpkin carriage, the crystal slipper, and the magic wand. The text \"Cinderella''s Ball\" should be featured prominently in the image, in fun and playful font.'\n\nfor ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $oldTitle) {\n        $p.Range.Delete()\n    }\n    elseif ($t -eq $oldDescription) {\n        $descRange = $p.Range\n        $descRange.Find.ClearFormatting()\n        [void]$descRange.Find.Execute($oldDescription)\n        $descRange.Text = $newPrompt\n    }\n}\n"}
